$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "247.62"
Set-TextValue "D3" "22.59"
Set-TextValue "D4" "5.286"
Set-TextValue "D5" "0.05729"
Set-TextValue "D7" "0.8088"
Set-TextValue "D8" "0.8692"
Set-TextValue "B9" "WazirX"
Set-TextValue "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1432"
Set-TextValue "E9" "8WazirXWRX"
Set-TextValue "B10" "MandalaExchangeToken"
Set-TextValue "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07368"
Set-TextValue "E10" "9MandalaExchangeTokenMDX"
Set-TextValue "B11" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03049"
Set-TextValue "E11" "10LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03115"
Set-TextValue "E12" "11BitrueCoinBTR"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09397"
Set-TextValue "E13" "12BitMartTokenBMX"
Set-TextValue "B14" "MCDex"
Set-TextValue "C14" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "3.877"
Set-TextValue "E14" "13MCDexMCB"
Set-TextValue "B15" "BitForexToken"
Set-TextValue "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001584"
Set-TextValue "E15" "14BitForexTokenBF"
Set-TextValue "B16" "CoinExToken"
Set-TextValue "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04817"
Set-TextValue "E16" "15CoinExTokenCET"
Set-TextValue "B17" "One"
Set-TextValue "C17" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005850"
Set-TextValue "E17" "16OneONE"
Set-TextValue "D18" "0.006145"
Set-TextValue "D19" "0.005161"
Set-TextValue "D20" "0.0009966"
Set-TextValue "D22" "3.727"
Set-TextValue "D23" "6.325"
Set-TextValue "D26" "0.1349"
Set-TextValue "D41" "0.006747"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.002931"
Set-TextValue "D44" "0.007298"
Set-TextValue "D45" "0.00005620"
Set-TextValue "D47" "0.6000"
Set-TextValue "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue "D48" "0.1774"
